$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text shared strings; runs share identical formatting) ---
$ws.Range("A8").Value = "Volume 31   Number  4"
$ws.Range("C9").Value = "Report Covering the Week  1/22/2024  Through  1/28/2024"

# --- Cells changing type: numeric -> shared-string "0" (keep style s=14) ---
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "0"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "0"

# Restore the original numeric-cell style (14) that got perturbed by the text-format switch above
$ws.Range("D15").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C23").PasteSpecial(-4122)

# --- Cells changing type: shared-string -> numeric (need style promoted to match numeric siblings) ---
$ws.Range("C20").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D23").PasteSpecial(-4122)

$ws.Range("H20").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("M20").PasteSpecial(-4122)
$ws.Range("M22").PasteSpecial(-4122)

$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("M20").Value = 500
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -100
$ws.Range("M22").Value = 0
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = -100

# --- Remaining straightforward numeric value updates (style/type unchanged) ---
$ws.Range("F15").Value = 1
$ws.Range("N15").Value = -66.666666666666
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("I16").Value = 10
$ws.Range("J16").Value = 6
$ws.Range("K16").Value = 66.666666666666
$ws.Range("L16").Value = -16.666666666666
$ws.Range("M16").Value = 25
$ws.Range("N16").Value = -78.260869565217
$ws.Range("E17").Value = -100
$ws.Range("F17").Value = 3
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = -70
$ws.Range("J17").Value = 10
$ws.Range("K17").Value = -70
$ws.Range("L17").Value = -57.142857142857
$ws.Range("M17").Value = -70
$ws.Range("N17").Value = -85.714285714285
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -66.666666666666
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 11
$ws.Range("J18").Value = 11
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 266.666666666667
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = -62.068965517241
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = -40
$ws.Range("F19").Value = 38
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = -20.833333333333
$ws.Range("I19").Value = 38
$ws.Range("J19").Value = 48
$ws.Range("K19").Value = -20.833333333333
$ws.Range("L19").Value = -26.923076923076
$ws.Range("M19").Value = -28.301886792452
$ws.Range("N19").Value = -29.629629629629
$ws.Range("C20").Value = 2
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 200
$ws.Range("I20").Value = 6
$ws.Range("J20").Value = 2
$ws.Range("K20").Value = 200
$ws.Range("L20").Value = -40
$ws.Range("N20").Value = -85.365853658536
$ws.Range("C21").Value = 14
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = -48.148148148148
$ws.Range("F21").Value = 69
$ws.Range("G21").Value = 77
$ws.Range("H21").Value = -10.38961038961
$ws.Range("I21").Value = 69
$ws.Range("J21").Value = 77
$ws.Range("K21").Value = -10.38961038961
$ws.Range("L21").Value = -17.857142857142
$ws.Range("M21").Value = -16.867469879518
$ws.Range("N21").Value = -64.432989690721
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -66.666666666666
$ws.Range("J22").Value = 3
$ws.Range("K22").Value = -66.666666666666
$ws.Range("L22").Value = -75
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = -33.333333333333
$ws.Range("J23").Value = 3
$ws.Range("K23").Value = -33.333333333333
$ws.Range("L23").Value = -71.428571428571
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 10
$ws.Range("E24").Value = 20
$ws.Range("F24").Value = 49
$ws.Range("G24").Value = 47
$ws.Range("H24").Value = 4.255319148936
$ws.Range("I24").Value = 49
$ws.Range("J24").Value = 47
$ws.Range("K24").Value = 4.255319148936
$ws.Range("L24").Value = 13.953488372093
$ws.Range("M24").Value = -20.967741935483
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 27
$ws.Range("G25").Value = 24
$ws.Range("H25").Value = 12.5
$ws.Range("I25").Value = 27
$ws.Range("J25").Value = 24
$ws.Range("K25").Value = 12.5
$ws.Range("L25").Value = 125
$ws.Range("F26").Value = 1
$ws.Range("C27").Value = 2
